$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching the formatting of the other
# header cells (copy formats from G1 "sum" then set the text).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add corresponding data value in H2
$ws.Range("H2").Value = 0
